# Apply odds updates to Sheet1 per the 2026-01-12 Betfair data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — Saudi Professional League: Al-Hazm (KSA) vs Al Najma Club
$ws.Range("G2").Value = 2.04
$ws.Range("I2").Value = 5
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.89
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 1.8
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.27
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 17.5
$ws.Range("Y2").Value = 19.5
$ws.Range("Z2").Value = 42
$ws.Range("AA2").Value = 130
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 22
$ws.Range("AE2").Value = 75
$ws.Range("AF2").Value = 14.5
$ws.Range("AG2").Value = 12.5
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 80
$ws.Range("AJ2").Value = 27
$ws.Range("AK2").Value = 25
$ws.Range("AL2").Value = 46
$ws.Range("AM2").Value = 130
$ws.Range("AN2").Value = 17
$ws.Range("AO2").Value = 85

# Row 3 — Cypriot 1st Division: Digenis Ypsona vs Olympiakos Nicosia FC
$ws.Range("Q3").Value = 2.14

# Row 4 — Cypriot 1st Division: Pafos FC vs Omonia FC Aradippou
$ws.Range("G4").Value = 1.22
$ws.Range("Q4").Value = 1.69

# Row 6 — Italian Serie A: Genoa vs Cagliari
$ws.Range("Z6").Value = 34

# Row 7 — Saudi Professional League: Al-Ettifaq vs Al-Khaleej Saihat
$ws.Range("G7").Value = 2.08
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 4.4
$ws.Range("P7").Value = 2.44
$ws.Range("Q7").Value = 1.55

# Row 8 — Saudi Professional League: Al-Hilal vs Al Nassr
$ws.Range("F8").Value = 2.28
$ws.Range("H8").Value = 2.64
$ws.Range("Q8").Value = 1.37

# Row 9 — Italian Serie C: Salernitana vs Cosenza
$ws.Range("F9").Value = 2.04
$ws.Range("I9").Value = 4.6
$ws.Range("J9").Value = 3.05
$ws.Range("Q9").Value = 1.99
$ws.Range("V9").Value = 1.28
$ws.Range("AD9").Value = 21

# Row 10 — Spanish Segunda Division: Huesca vs Cordoba
$ws.Range("Q10").Value = 2.36

# Row 11 — Italian Serie A: Juventus vs US Cremonese
$ws.Range("G11").Value = 1.29
$ws.Range("R11").Value = 1.54
$ws.Range("Z11").Value = 160

# Row 12 — Spanish La Liga: Sevilla vs Celta Vigo
$ws.Range("F12").Value = 2.84
$ws.Range("I12").Value = 2.86
$ws.Range("K12").Value = 3.4
